# "All Modules work" commit:
#  - Student roster names were updated on the "فیزیک" (Physics) and
#    "شیمی" (Chemistry) sheets (row 4 / "علی شمس" stays the same,
#    every other student A2:A13 got a new name). The "ریاضی" (Math)
#    sheet's roster is untouched.
#  - The Math sheet ("ریاضی") ends up the active sheet/tab, with the
#    cursor resting on B16 (below the data table).

$wb = $excel.ActiveWorkbook

$wsMath    = $wb.Worksheets.Item(1)   # ریاضی
$wsPhysics = $wb.Worksheets.Item(2)   # فیزیک
$wsChem    = $wb.Worksheets.Item(3)   # شیمی

# New student roster for Physics & Chemistry (row 4 "علی شمس" unchanged).
$newNames = @(
    "رضا مولایی",
    "محمود وکیلی",
    "علی شمس",
    "رضا برهانی مرند",
    "سینا وکیلی",
    "رضا خوشخو",
    "کسری نیک فرجام",
    "سید محسن ابطحی",
    "ندا اشرفی",
    "اکرم سینایی",
    "پوریا مقدسی",
    "احمد رضا معین"
)

for ($i = 0; $i -lt $newNames.Count; $i++) {
    $row = $i + 2
    $wsPhysics.Range("A$row").Value = $newNames[$i]
    $wsChem.Range("A$row").Value = $newNames[$i]
}

# Math sheet becomes the active tab with the cursor on B16.
$wsMath.Activate() | Out-Null
$wsMath.Range("B16").Select() | Out-Null
